# "Add files via upload" / "Updated data tables."
#
# Appends one new datastore record (id=16, a non-numeric-formatted variant
# of the "Budget vs Actual" table-valued function) to the
# analytics_datastores sheet, and the matching parameter row (Amp/VARCHAR)
# to the analytics_datastores_parameters sheet.

$wb = $excel.ActiveWorkbook

$wsDatastores  = $wb.Worksheets.Item("analytics_datastores")
$wsParameters  = $wb.Worksheets.Item("analytics_datastores_parameters")

# --- analytics_datastores: new row 21 (id 16) -------------------------------
$row = 21
$wsDatastores.Cells.Item($row, 1).Value = 16
$wsDatastores.Cells.Item($row, 2).Value = "B"
$wsDatastores.Cells.Item($row, 3).Value = "Budget Reports"
$wsDatastores.Cells.Item($row, 4).Value = "Budget vs Actual"
$wsDatastores.Cells.Item($row, 5).Value = 3
$wsDatastores.Cells.Item($row, 6).Value = "Budget vs Actual"
$wsDatastores.Cells.Item($row, 7).Value = "Same as id=9, except does not format the numerical columns. This is for Excel, where using id=9 renders text formats that don't render correctly."
$wsDatastores.Cells.Item($row, 8).Value = "ITAnalytics"
$wsDatastores.Cells.Item($row, 9).Value = "SELECT * FROM [ITAnalytics].[dbo].[fn_BudgetActual_no_numeric_format](Amp) ORDER BY ACCT"
$wsDatastores.Cells.Item($row, 10).Value = 1

# --- analytics_datastores_parameters: new row 8 -----------------------------
$prow = 8
$wsParameters.Cells.Item($prow, 1).Value = 3
$wsParameters.Cells.Item($prow, 2).Value = 16
$wsParameters.Cells.Item($prow, 3).Value = "Amp"
$wsParameters.Cells.Item($prow, 4).Value = "VARCHAR"
$wsParameters.Cells.Item($prow, 5).Value = 6
